$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.249960064888
$ws.Range("B1").Value = 2.038832426071167
$ws.Range("C1").Value = 5.871021747589111
$ws.Range("D1").Value = 1.956021070480347
$ws.Range("E1").Value = 1.13567328453064
